# Generate Report for Handback
#
# For each locale sheet (zh-cn, de-de):
#   - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#   - A "Latest Target File" (F) / "Latest Handback File" (G) pair of hyperlinked
#     cells is populated for every data row, mirroring the existing
#     Source File (A) / Latest Handoff File (D) hyperlinks
#   - "Latest Handback DateTime" (H) is stamped with the real handback time
#     (previously the zero-date placeholder)

$wb = $excel.ActiveWorkbook
$statusText = "Handed back: in sync with en-US"

function Find-HyperlinkAtAddress($ws, $addr) {
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq $addr) {
            return $hl
        }
    }
    return $null
}

function Remove-HyperlinkAtAddress($ws, $addr) {
    $hl = Find-HyperlinkAtAddress $ws $addr
    if ($hl) { $hl.Delete() }
}

function Set-HandbackRow($ws, $row, $handbackStamp) {
    $aAddr = $ws.Cells.Item($row, 1).Address()
    $dAddr = $ws.Cells.Item($row, 4).Address()

    # Status column
    $ws.Cells.Item($row, 3).Value = $statusText

    # Mirror the source (A) / handoff (D) hyperlinks into the new
    # target (F) / handback (G) columns.
    $srcLink = Find-HyperlinkAtAddress $ws $aAddr
    $offLink = Find-HyperlinkAtAddress $ws $dAddr

    $ws.Hyperlinks.Add($ws.Cells.Item($row, 6), $srcLink.Address, [Type]::Missing, [Type]::Missing, $srcLink.TextToDisplay) | Out-Null
    $ws.Hyperlinks.Add($ws.Cells.Item($row, 7), $offLink.Address, [Type]::Missing, [Type]::Missing, $offLink.TextToDisplay) | Out-Null

    # Latest Handback DateTime
    $ws.Cells.Item($row, 8).Value = $handbackStamp
}

function Update-LocaleSheet($ws, $handbackStamp) {
    # Row 3's Source (A3) / Handoff (D3) hyperlinks are re-created (rather
    # than left alone) so the relationship ids line up after row 2 gains its
    # two new link columns.
    $a3Addr = $ws.Cells.Item(3, 1).Address()
    $d3Addr = $ws.Cells.Item(3, 4).Address()

    $row3Src = Find-HyperlinkAtAddress $ws $a3Addr
    $row3SrcAddress = $row3Src.Address
    $row3SrcDisplay = $row3Src.TextToDisplay

    $row3Off = Find-HyperlinkAtAddress $ws $d3Addr
    $row3OffAddress = $row3Off.Address
    $row3OffDisplay = $row3Off.TextToDisplay

    Remove-HyperlinkAtAddress $ws $a3Addr
    Remove-HyperlinkAtAddress $ws $d3Addr

    Set-HandbackRow $ws 2 $handbackStamp

    $ws.Hyperlinks.Add($ws.Range("A3"), $row3SrcAddress, [Type]::Missing, [Type]::Missing, $row3SrcDisplay) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("D3"), $row3OffAddress, [Type]::Missing, [Type]::Missing, $row3OffDisplay) | Out-Null

    Set-HandbackRow $ws 3 $handbackStamp
}

$wsZh = $wb.Worksheets.Item("zh-cn")
Update-LocaleSheet $wsZh "2016-03-31 07:53:42"

$wsDe = $wb.Worksheets.Item("de-de")
Update-LocaleSheet $wsDe "2016-03-31 07:53:59"

Write-Host "Handback report generated."
